$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added to the daily log. They sort into
# the existing date-ordered block at rows 258-259, pushing the previous
# contents of rows 258-285 down to rows 260-287.
$ws.Rows("258:259").Insert()

# Row 258: new record
$ws.Range("A258").Value = 6
$ws.Range("B258").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C258").Value = "Metropolitana"
$ws.Range("D258").Value = [datetime]"2022-01-17"
$ws.Range("E258").Value = 13
$ws.Range("F258").Value = "Fruta"
$ws.Range("G258").Value = 100101
$ws.Range("H258").Value = "Berries"
$ws.Range("I258").Value = 100101001
$ws.Range("J258").Value = "Arándano (blue)"
$ws.Range("K258").Value = "Sin especificar"
$ws.Range("L258").Value = "Primera"
$ws.Range("M258").Value = 250
$ws.Range("N258").Value = 4000
$ws.Range("O258").Value = 4000
$ws.Range("P258").Value = 4000
$ws.Range("Q258").Value = "$/bandeja 2 kilos"
$ws.Range("R258").Value = "Provincia de Curicó"
$ws.Range("S258").Value = 2000
$ws.Range("T258").Value = 2

# Row 259: new record
$ws.Range("A259").Value = 6
$ws.Range("B259").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C259").Value = "Metropolitana"
$ws.Range("D259").Value = [datetime]"2022-01-17"
$ws.Range("E259").Value = 13
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100101
$ws.Range("H259").Value = "Berries"
$ws.Range("I259").Value = 100101001
$ws.Range("J259").Value = "Arándano (blue)"
$ws.Range("K259").Value = "Sin especificar"
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value = 300
$ws.Range("N259").Value = 3400
$ws.Range("O259").Value = 3400
$ws.Range("P259").Value = 3400
$ws.Range("Q259").Value = "$/bandeja 2 kilos"
$ws.Range("R259").Value = "Provincia de Linares"
$ws.Range("S259").Value = 1700
$ws.Range("T259").Value = 2
